$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.959.22"
$ws.Range("E2").Value = "  -0.67%  "

$ws.Range("D3").Value = "1.743.79"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "248.36"
$ws.Range("E5").Value = "  +5.07%  "

$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").Value = "0.5053"
$ws.Range("E7").Value = "  -4.51%  "

$ws.Range("D8").Value = "0.2743"
$ws.Range("E8").Value = "  -2.93%  "

$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("E10").Value = "  +1.11%  "

$ws.Range("D11").Value = "1.742.01"
$ws.Range("E11").Value = "  -0.37%  "

$ws.Range("D12").Value = "0.6542"
$ws.Range("E12").Value = "  +1.28%  "

$ws.Range("E13").Value = "  -1.76%  "

$ws.Range("D14").Value = "4.643"
$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("E15").Value = "  -1.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9990"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").Value = "0.9988"
$ws.Range("E17").Value = "  -0.07%  "

$ws.Range("D18").Value = "25.986.49"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19").Value = "11.84"

$ws.Range("D20").Value = "0.000006838"
$ws.Range("E20").Value = "  +1.39%  "

$ws.Range("D21").Value = "1.969.13"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").Value = "4.385"
$ws.Range("E22").Value = "  +1.32%  "

$ws.Range("D23").Value = "8.694"
$ws.Range("E23").Value = "  -0.27%  "

$ws.Range("D24").Value = "5.403"
$ws.Range("E24").Value = "  +3.33%  "

$ws.Range("D25").Value = "136.66"
$ws.Range("E25").Value = "  -2.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.500"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.35%  "

$ws.Range("E27").Value = "  -0.29%  "

$ws.Range("D28").Value = "1.779"
$ws.Range("E28").Value = "  -1.40%  "

$ws.Range("D29").Value = "105.46"
$ws.Range("E29").Value = "  +0.60%  "

$ws.Range("D30").Value = "3.897"
$ws.Range("E30").Value = "  +2.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08240"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.86%  "

$ws.Range("D32").Value = "3.638"
$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04680"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.95%  "

$ws.Range("D34").Value = "2.655"
$ws.Range("E34").Value = "  +0.40%  "

$ws.Range("D35").Value = "0.9932"
$ws.Range("E35").Value = "  -1.45%  "

$ws.Range("D36").Value = "0.6187"
$ws.Range("E36").Value = "  -2.51%  "

$ws.Range("D37").Value = "2.752"
$ws.Range("E37").Value = "  +1.75%  "

$ws.Range("E38").Value = "  -1.09%  "

$ws.Range("D39").Value = "1.928"
$ws.Range("E39").Value = "  -2.30%  "

$ws.Range("D40").Value = "0.9989"

$ws.Range("D41").Value = "100.07"
$ws.Range("E41").Value = "  -2.45%  "

$ws.Range("D42").Value = "0.3913"
$ws.Range("E42").Value = "  -0.39%  "

$ws.Range("D43").Value = "0.7589"
$ws.Range("E43").Value = "  +1.20%  "

$ws.Range("D44").Value = "5.013"
$ws.Range("E44").Value = "  -0.80%  "

$ws.Range("D45").Value = "0.1146"
$ws.Range("E45").Value = "  -0.60%  "

$ws.Range("D46").Value = "6.296"
$ws.Range("E46").Value = "  -0.89%  "

$ws.Range("D47").Value = "55.59"
$ws.Range("E47").Value = "  +1.90%  "

$ws.Range("E48").Value = "  -1.58%  "

$ws.Range("D49").Value = "30.61"
$ws.Range("E49").Value = "  -1.34%  "

$ws.Range("D50").Value = "7.562"
$ws.Range("E50").Value = "  -0.40%  "

$ws.Range("D51").Value = "0.3435"
$ws.Range("E51").Value = "  -1.19%  "
